# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts). Update the values for rows 2-14 to match
# the regenerated save data.
$kValues = @{
    2  = 4
    3  = 4
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 5
    9  = 5
    10 = 1
    11 = 1
    12 = 4
    13 = 1
    14 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
